$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = '[''Argentina'', ''Italy'', ''Bulgaria'']'
$ws.Range("H3").Value = '[''Mexico'', ''Paraguay'', ''Belgium'']'
$ws.Range("H5").Value = '[''Northern Ireland'', ''Brazil'', ''Spain'']'
$ws.Range("H6").Value = '[''Algeria'', ''Brazil'', ''Spain'']'
$ws.Range("H7").Value = '[''Northern Ireland'', ''Brazil'', ''Spain'']'
$ws.Range("H8").Value = '[''Uruguay'', ''Denmark'', ''West Germany'']'
$ws.Range("H9").Value = '[''Portugal'', ''Morocco'', ''Poland'']'
$ws.Range("H10").Value = '[''Portugal'', ''England'', ''Poland'']'
$ws.Range("H11").Value = '[''Portugal'', ''Morocco'', ''England'']'
$ws.Range("H12").Value = '[''Morocco'', ''England'', ''Poland'']'
$ws.Range("H13").Value = '[''Italy'', ''Czechoslovakia'', ''Austria'']'
$ws.Range("H14").Value = '[''Argentina'', ''Romania'', ''Cameroon'']'
$ws.Range("H15").Value = '[''Argentina'', ''Soviet Union'', ''Cameroon'']'
$ws.Range("H16").Value = '[''Argentina'', ''Romania'', ''Cameroon'']'
$ws.Range("H17").Value = '[''Costa Rica'', ''Scotland'', ''Brazil'']'
$ws.Range("H19").Value = '[''Costa Rica'', ''Scotland'', ''Brazil'']'
$ws.Range("H20").Value = '[''Yugoslavia'', ''Colombia'', ''West Germany'']'
$ws.Range("H21").Value = '[''Belgium'', ''Uruguay'', ''Spain'']'
$ws.Range("H22").Value = '[''Netherlands'', ''England'', ''Republic of Ireland'']'
$ws.Range("H23").Value = '[''Egypt'', ''Netherlands'', ''England'']'
$ws.Range("H24").Value = '[''Netherlands'', ''England'', ''Republic of Ireland'']'
$ws.Range("H25").Value = '[''United States'', ''Switzerland'', ''Romania'']'
$ws.Range("H27").Value = '[''Sweden'', ''Russia'', ''Brazil'']'
$ws.Range("H28").Value = '[''Germany'', ''South Korea'', ''Spain'']'
$ws.Range("H29").Value = '[''Argentina'', ''Nigeria'', ''Bulgaria'']'
$ws.Range("H30").Value = '[''Mexico'', ''Italy'', ''Republic of Ireland'']'
$ws.Range("H31").Value = '[''Norway'', ''Italy'', ''Republic of Ireland'']'
$ws.Range("H32").Value = '[''Mexico'', ''Italy'', ''Republic of Ireland'']'
$ws.Range("H35").Value = '[''Morocco'', ''Brazil'']'
$ws.Range("H37").Value = '[''Chile'', ''Italy'']'
$ws.Range("H40").Value = '[''Nigeria'', ''Spain'']'
$ws.Range("H42").Value = '[''Netherlands'', ''Mexico'']'
$ws.Range("H43").Value = '[''Netherlands'', ''Belgium'']'
$ws.Range("H44").Value = '[''Netherlands'', ''Mexico'']'
$ws.Range("H47").Value = '[''Argentina'', ''Croatia'']'
$ws.Range("H48").Value = '[''Senegal'', ''Denmark'']'
$ws.Range("H50").Value = '[''Paraguay'', ''Spain'']'
$ws.Range("H51").Value = '[''Costa Rica'', ''Brazil'']'
$ws.Range("H52").Value = '[''Turkey'', ''Brazil'']'
$ws.Range("H53").Value = '[''Costa Rica'', ''Brazil'']'
$ws.Range("H54").Value = '[''Turkey'', ''Brazil'']'
$ws.Range("H55").Value = '[''United States'', ''South Korea'']'
$ws.Range("H57").Value = '[''United States'', ''South Korea'']'
$ws.Range("H59").Value = '[''Germany'', ''Republic of Ireland'']'
$ws.Range("H64").Value = '[''Russia'', ''Japan'']'
$ws.Range("H65").Value = '[''Belgium'', ''Japan'']'
$ws.Range("H66").Value = '[''Russia'', ''Japan'']'
$ws.Range("H67").Value = '[''Belgium'', ''Japan'']'
$ws.Range("H70").Value = '[''Argentina'', ''Netherlands'']'
$ws.Range("H72").Value = '[''Czech Republic'', ''Italy'']'
$ws.Range("H79").Value = '[''Switzerland'', ''South Korea'']'
$ws.Range("H80").Value = '[''Switzerland'', ''France'']'
$ws.Range("H83").Value = '[''Argentina'', ''South Korea'']'
$ws.Range("H84").Value = '[''Argentina'', ''Greece'']'
$ws.Range("H85").Value = '[''Argentina'', ''South Korea'']'
$ws.Range("H86").Value = '[''United States'', ''Slovenia'']'
$ws.Range("H87").Value = '[''England'', ''Slovenia'']'
$ws.Range("H88").Value = '[''United States'', ''England'']'
$ws.Range("H89").Value = '[''Germany'', ''Ghana'']'
$ws.Range("H90").Value = '[''Netherlands'', ''Japan'']'
$ws.Range("H94").Value = '[''Chile'', ''Spain'']'
$ws.Range("H97").Value = '[''Ivory Coast'', ''Colombia'']'
$ws.Range("H98").Value = '[''Greece'', ''Colombia'']'
$ws.Range("H99").Value = '[''Ivory Coast'', ''Colombia'']'
$ws.Range("H100").Value = '[''Greece'', ''Colombia'']'
$ws.Range("H101").Value = '[''Costa Rica'', ''Italy'']'
$ws.Range("H102").Value = '[''Costa Rica'', ''Uruguay'']'
$ws.Range("H103").Value = '[''Ecuador'', ''France'']'
$ws.Range("H104").Value = '[''Switzerland'', ''France'']'
$ws.Range("H105").Value = '[''Argentina'', ''Nigeria'']'
$ws.Range("H107").Value = '[''Belgium'', ''Algeria'']'
$ws.Range("H109").Value = '[''Belgium'', ''Algeria'']'
$ws.Range("H113").Value = '[''Nigeria'', ''Croatia'']'
$ws.Range("H114").Value = '[''Argentina'', ''Croatia'']'
$ws.Range("H115").Value = '[''Nigeria'', ''Croatia'']'
$ws.Range("H116").Value = '[''Argentina'', ''Croatia'']'
$ws.Range("H117").Value = '[''Switzerland'', ''Brazil'']'
$ws.Range("H121").Value = '[''Senegal'', ''Japan'']'
$ws.Range("H124").Value = '[''Netherlands'', ''Senegal'']'
$ws.Range("H126").Value = '[''Netherlands'', ''Senegal'']'
$ws.Range("H128").Value = '[''United States'', ''England'']'
$ws.Range("H129").Value = '[''Argentina'', ''Poland'']'
$ws.Range("H134").Value = '[''Germany'', ''Spain'']'
$ws.Range("H135").Value = '[''Spain'', ''Japan'']'
$ws.Range("H136").Value = '[''Costa Rica'', ''Japan'']'
$ws.Range("H137").Value = '[''Spain'', ''Japan'']'
$ws.Range("H138").Value = '[''Morocco'', ''Croatia'']'
$ws.Range("H139").Value = '[''Switzerland'', ''Brazil'']'
$ws.Range("H140").Value = '[''Brazil'', ''Serbia'']'
$ws.Range("H141").Value = '[''Switzerland'', ''Brazil'']'
